# Slide 1, title placeholder shape: merge the "QUALITY " and "MONITORING"
# runs into a single run "QUALITY MONITORING" (as a real user would do by
# selecting across the two runs and retyping them as one continuous run).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Full title text is: "WATER QUALITY MONITORING<break>SYSTEM  USING IOT"
# Characters 7-24 (1-based, length 18) cover exactly "QUALITY " + "MONITORING".
$target = $tr.Characters(7, 18)
$target.Text = "QUALITY MONITORING"
